# Apply updated "dSF" (column F) values for several rows in Sheet1.
# This reflects a repull/recalculation of data (mean calculation) where
# only the dSF column values changed for the listed rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 1
    6  = -5
    8  = -4
    12 = -2
    15 = -6
    16 = 6
    17 = -6
    18 = -5
    19 = -8
    21 = -4
    23 = -6
    31 = -3
    33 = 3
    35 = -3
    37 = 1
    39 = 5
    43 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
